# Update recomputed lower/upper confidence-interval bound values
# (columns F="lower", G="upper") across the RLIe findings sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 0.8
$ws.Range("G14").Value = 0.7067669172932332
$ws.Range("F15").Value = 0.5729323308270677
$ws.Range("G15").Value = 0.6977443609022556
$ws.Range("G16").Value = 0.6962406015037594
$ws.Range("F17").Value = 0.5684210526315789
$ws.Range("G17").Value = 0.6947368421052631
$ws.Range("G19").Value = 0.8844444444444445
$ws.Range("F20").Value = 0.7999722222222223
$ws.Range("G20").Value = 0.8766666666666667
$ws.Range("F21").Value = 0.7933333333333333
$ws.Range("F22").Value = 0.8934239130434783
$ws.Range("F24").Value = 0.8673913043478261
$ws.Range("G24").Value = 0.9630434782608696
$ws.Range("F25").Value = 0.8695652173913043
$ws.Range("G25").Value = 0.9630434782608696
$ws.Range("G29").Value = 0.78
$ws.Range("F42").Value = 0.5833333333333333
$ws.Range("G42").Value = 0.7060606060606061
$ws.Range("G43").Value = 0.6909090909090909
$ws.Range("G44").Value = 0.6909090909090909
$ws.Range("F45").Value = 0.5696969696969697
$ws.Range("G45").Value = 0.6909090909090909
$ws.Range("F46").Value = 0.8842105263157894
$ws.Range("F47").Value = 0.8789473684210526
$ws.Range("F48").Value = 0.8631578947368421
$ws.Range("F49").Value = 0.8631578947368421
$ws.Range("F51").Value = 0.7462686567164178
$ws.Range("G52").Value = 0.8746268656716418
$ws.Range("F53").Value = 0.7373134328358208
$ws.Range("G53").Value = 0.8746268656716418
$ws.Range("F55").Value = 0.76
$ws.Range("G55").Value = 0.888
$ws.Range("F56").Value = 0.7573333333333333
$ws.Range("G56").Value = 0.8853333333333333
$ws.Range("F57").Value = 0.7493333333333334
$ws.Range("G57").Value = 0.8773333333333333
$ws.Range("F62").Value = 0.8603174603174604
$ws.Range("F64").Value = 0.8444444444444444
$ws.Range("F65").Value = 0.8412698412698413
$ws.Range("F70").Value = 0.7943844492440605
$ws.Range("G70").Value = 0.8479481641468682
$ws.Range("F71").Value = 0.7697624190064795
$ws.Range("F72").Value = 0.767170626349892
$ws.Range("G72").Value = 0.8241900647948164
$ws.Range("G73").Value = 0.8220302375809936
